$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row in column C (the "Förändrad" / Changed date column)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = $ws.UsedRange.Rows.Count }

# All data rows in column C currently hold the old "last changed" date serial
# value (45177 = 2023-09-08). Bump them all to the new value (45178 = 2023-09-09).
$ws.Range("C2:C$lastRow").Value = 45178
